$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.920350755611925
$ws.Range("D2").Value = 4.135355856447947
$ws.Range("E2").Value = 16.52162517355821
$ws.Range("F2").Value = 23.18408403552703
$ws.Range("G2").Value = 3.609790598073208
$ws.Range("K2").Value = 13.72358049572382
$ws.Range("N2").Value = 16.29667848184512
$ws.Range("O2").Value = 20.47193895262672
$ws.Range("C3").Value = 4.749567624056075
$ws.Range("D3").Value = 4.121053021641644
$ws.Range("E3").Value = 15.58156061696597
$ws.Range("F3").Value = 23.08901528907493
$ws.Range("G3").Value = 3.612358685666941
$ws.Range("K3").Value = 13.14297119694198
$ws.Range("N3").Value = 16.34121424621785
$ws.Range("O3").Value = 20.47306082942276
$ws.Range("C4").Value = 4.643060145993967
$ws.Range("D4").Value = 4.11223707160597
$ws.Range("E4").Value = 14.97957795359731
$ws.Range("F4").Value = 23.03930896789242
$ws.Range("G4").Value = 3.614016723684896
$ws.Range("K4").Value = 12.77466503047325
$ws.Range("N4").Value = 16.37042260430349
$ws.Range("O4").Value = 20.48076765352159
$ws.Range("C5").Value = 4.599321554243794
$ws.Range("D5").Value = 4.108637312335661
$ws.Range("E5").Value = 14.72830362133381
$ws.Range("F5").Value = 23.02124448453047
$ws.Range("G5").Value = 3.61471288497807
$ws.Range("K5").Value = 12.62182257761287
$ws.Range("N5").Value = 16.38279459472487
$ws.Range("O5").Value = 20.4856665469432
$ws.Range("C6").Value = 4.592040904730997
$ws.Range("D6").Value = 4.108039180299733
$ws.Range("E6").Value = 14.6862285725776
$ws.Range("F6").Value = 23.01837751562321
$ws.Range("G6").Value = 3.614829722175771
$ws.Range("K6").Value = 12.59628365029895
$ws.Range("N6").Value = 16.38487732642147
$ws.Range("O6").Value = 20.48658598499536
$ws.Range("C7").Value = 4.642471524612306
$ws.Range("D7").Value = 4.112188551345091
$ws.Range("E7").Value = 14.97621292286132
$ws.Range("F7").Value = 23.03905645864957
$ws.Range("G7").Value = 3.614026029269456
$ws.Range("K7").Value = 12.77261459319311
$ws.Range("N7").Value = 16.37058755584182
$ws.Range("O7").Value = 20.48082661211014
$ws.Range("C8").Value = 4.861851176889172
$ws.Range("D8").Value = 4.130432093457166
$ws.Range("E8").Value = 16.20277977941314
$ws.Range("F8").Value = 23.14951399672527
$ws.Range("G8").Value = 3.610659261391275
$ws.Range("K8").Value = 13.52595812019063
$ws.Range("N8").Value = 16.31164819334096
$ws.Range("O8").Value = 20.47086576328111
$ws.Range("C9").Value = 5.27589143431292
$ws.Range("D9").Value = 4.165879818106048
$ws.Range("E9").Value = 18.49744056464214
$ws.Range("F9").Value = 23.43417899362209
$ws.Range("G9").Value = 3.604698142474206
$ws.Range("K9").Value = 14.9014477123289
$ws.Range("N9").Value = 16.21081597262286
$ws.Range("O9").Value = 20.50723730979496
$ws.Range("C10").Value = 5.566533790561478
$ws.Range("D10").Value = 4.191646386625841
$ws.Range("E10").Value = 20.14641657501729
$ws.Range("F10").Value = 23.68362224432179
$ws.Range("G10").Value = 3.600704635920741
$ws.Range("K10").Value = 15.84078121051698
$ws.Range("N10").Value = 16.14567682025541
$ws.Range("O10").Value = 20.56824215504414
$ws.Range("C11").Value = 5.695156475909801
$ws.Range("D11").Value = 4.20329184352657
$ws.Range("E11").Value = 20.85447156346804
$ws.Range("F11").Value = 23.80552594241502
$ws.Range("G11").Value = 3.598970725321095
$ws.Range("K11").Value = 16.25113018631928
$ws.Range("N11").Value = 16.11797533431665
$ws.Range("O11").Value = 20.60344869814123
$ws.Range("C12").Value = 5.743299510532393
$ws.Range("D12").Value = 4.207689364995182
$ws.Range("E12").Value = 21.11658107441131
$ws.Range("F12").Value = 23.852867898235
$ws.Range("G12").Value = 3.598325961821081
$ws.Range("K12").Value = 16.40397850135355
$ws.Range("N12").Value = 16.10776242895233
$ws.Range("O12").Value = 20.61785057603461
$ws.Range("C13").Value = 5.732956802721037
$ws.Range("D13").Value = 4.206742852777809
$ws.Range("E13").Value = 21.06039809866114
$ws.Range("F13").Value = 23.84262004490043
$ws.Range("G13").Value = 3.598464298098452
$ws.Range("K13").Value = 16.37117441338799
$ws.Range("N13").Value = 16.10994964804346
$ws.Range("O13").Value = 20.61470134018847
$ws.Range("C14").Value = 5.699128777261166
$ws.Range("D14").Value = 4.203653886258087
$ws.Range("E14").Value = 20.87615583476171
$ws.Range("F14").Value = 23.80939734841208
$ws.Range("G14").Value = 3.598917443528873
$ws.Range("K14").Value = 16.26375658914797
$ws.Range("N14").Value = 16.11712956323889
$ws.Range("O14").Value = 20.60461212448306
$ws.Range("C15").Value = 5.678333437775838
$ws.Range("D15").Value = 4.201760152516348
$ws.Range("E15").Value = 20.76251974456583
$ws.Range("F15").Value = 23.7892000975781
$ws.Range("G15").Value = 3.599196546857907
$ws.Range("K15").Value = 16.19762623901013
$ws.Range("N15").Value = 16.12156353153102
$ws.Range("O15").Value = 20.59857143296649
$ws.Range("C16").Value = 5.558051564537267
$ws.Range("D16").Value = 4.190883677703427
$ws.Range("E16").Value = 20.09929985818864
$ws.Range("F16").Value = 23.67582247307809
$ws.Range("G16").Value = 3.600819610374577
$ws.Range("K16").Value = 15.81361349156376
$ws.Range("N16").Value = 16.14752597691723
$ws.Range("O16").Value = 20.56609119834584
$ws.Range("C17").Value = 5.483307258525618
$ws.Range("D17").Value = 4.184190844598014
$ws.Range("E17").Value = 19.68168473264206
$ws.Range("F17").Value = 23.60840597465388
$ws.Range("G17").Value = 3.601836453514804
$ws.Range("K17").Value = 15.57361193745514
$ws.Range("N17").Value = 16.16394715971218
$ws.Range("O17").Value = 20.54807409057572
$ws.Range("C18").Value = 5.439980659305506
$ws.Range("D18").Value = 4.180334247715995
$ws.Range("E18").Value = 19.43752145254371
$ws.Range("F18").Value = 23.57042547402282
$ws.Range("G18").Value = 3.602429108265574
$ws.Range("K18").Value = 15.43397825480801
$ws.Range("N18").Value = 16.17357394362401
$ws.Range("O18").Value = 20.5384130723129
$ws.Range("C19").Value = 5.425254895163467
$ws.Range("D19").Value = 4.179027302566878
$ws.Range("E19").Value = 19.3541704138418
$ws.Range("F19").Value = 23.55770353616082
$ws.Range("G19").Value = 3.602631111634169
$ws.Range("K19").Value = 15.38643082307079
$ws.Range("N19").Value = 16.17686464168683
$ws.Range("O19").Value = 20.53526261432931
$ws.Range("C20").Value = 5.491299014781005
$ws.Range("D20").Value = 4.18490404895698
$ws.Range("E20").Value = 19.7265503876805
$ws.Range("F20").Value = 23.6155004531315
$ws.Range("G20").Value = 3.601727402737732
$ws.Range("K20").Value = 15.5993259883476
$ws.Range("N20").Value = 16.16218029028191
$ws.Range("O20").Value = 20.54991940399786
$ws.Range("C21").Value = 5.709080531422051
$ws.Range("D21").Value = 4.204561537566619
$ws.Range("E21").Value = 20.93043524068451
$ws.Range("F21").Value = 23.81912391811318
$ws.Range("G21").Value = 3.598784023197424
$ws.Range("K21").Value = 16.29537754063675
$ws.Range("N21").Value = 16.11501313361146
$ws.Range("O21").Value = 20.60754655939306
$ws.Range("C22").Value = 5.848110341199933
$ws.Range("D22").Value = 4.217336130351911
$ws.Range("E22").Value = 21.68219620283627
$ws.Range("F22").Value = 23.95906255092576
$ws.Range("G22").Value = 3.59692928160653
$ws.Range("K22").Value = 16.7354339275301
$ws.Range("N22").Value = 16.08580125899585
$ws.Range("O22").Value = 20.65144341647732
$ws.Range("C23").Value = 5.774223872020617
$ws.Range("D23").Value = 4.21052521967101
$ws.Range("E23").Value = 21.28416274900769
$ws.Range("F23").Value = 23.88375858690698
$ws.Range("G23").Value = 3.597912908088974
$ws.Range("K23").Value = 16.50195679561365
$ws.Range("N23").Value = 16.10124463478961
$ws.Range("O23").Value = 20.62744553912198
$ws.Range("C24").Value = 5.487687044231699
$ws.Range("D24").Value = 4.184581636538844
$ws.Range("E24").Value = 19.70627930760945
$ws.Range("F24").Value = 23.61229061310861
$ws.Range("G24").Value = 3.601776679488712
$ws.Range("K24").Value = 15.58770580849445
$ws.Range("N24").Value = 16.16297851243684
$ws.Range("O24").Value = 20.54908296574273
$ws.Range("C25").Value = 5.16601981502649
$ws.Range("D25").Value = 4.156333317682551
$ws.Range("E25").Value = 17.85265245399885
$ws.Range("F25").Value = 23.34999252014823
$ws.Range("G25").Value = 3.606242636356533
$ws.Range("K25").Value = 14.54127465749413
$ws.Range("N25").Value = 16.2365200188154
$ws.Range("O25").Value = 20.49138572767788
